$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $old"
    }
    return $ok
}

# 1. "Is more likely..." [-4, -3] -> [-4, -2]
ReplaceText "Is more likely to contain the population value in the [-1,1] range than in the [-4, -3] range" "Is more likely to contain the population value in the [-1,1] range than in the [-4, -2] range"

# 2. "Is less likely..." [-4, -3] -> [-4, -2], then move _GoBack bookmark here (right after the new "2")
$rng2 = $d.Content
$rng2.Find.Execute("Is less likely to contain the population value in the [-1,1] range than in the [-4, -3] range", $true, $false, $false, $false, $false, $true, 1, $false, "Is less likely to contain the population value in the [-1,1] range than in the [-4, -2] range", 2)
$bmPos = $rng2.End - 7   # length of "] range" = 7, positions bookmark right after "...-2"
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. "(Similar to Section 7) " -> "(35 -37, Similar to Section 7) "
ReplaceText "(Similar to Section 7) " "(35 -37, Similar to Section 7) "

# 4. "Bring text from exchange between commentator and interviewer on gender. "
#    -> "Go back to lecture 1, and review the exchange between an interviewer and a commentator on the gender wage gap. "
ReplaceText "Bring text from exchange between commentator and interviewer on gender. " "Go back to lecture 1, and review the exchange between an interviewer and a commentator on the gender wage gap. "

# 5. "...women earn 9% less than men on average." -> "...women earn 9% less than men on average in the UK."
ReplaceText "Write down the underlying regression that corresponds to interviewers claim that women earn 9% less than men on average." "Write down the underlying regression that corresponds to interviewers claim that women earn 9% less than men on average in the UK."

# 6. "...corresponds to his response that …" -> "...corresponds to the commentator's response"
ReplaceText "Write down the underlying regression that corresponds to his response that …" "Write down the underlying regression that corresponds to the commentator’s response"

# 7. "...regarding the first estimate? " -> "...regarding the first estimate (in a)? "
ReplaceText "Using concepts from class, particularly of regression, what is the commentator implying regarding the first estimate? " "Using concepts from class, particularly of regression, what is the commentator implying regarding the first estimate (in a)? "

# 8. “experience” -> “personality”
ReplaceText "Discuss how the OVB formula could help us understand the effect of including a variable like “experience” on the effect of gender on wages. " "Discuss how the OVB formula could help us understand the effect of including a variable like “personality” on the effect of gender on wages. "

# 9. Insert new blank paragraph after "...as a motivation for learning the remaining material of the course!)"
$rng9 = $d.Content
$rng9.Find.Execute("as a motivation for learning the remaining material of the course!)")
$para9 = $rng9.Paragraphs(1)
$prange9 = $para9.Range
$prange9.Collapse(0)
$prange9.InsertParagraphAfter()

# 15/16. Remove old _GoBack bookmark implicitly handled by step 2 (Bookmarks.Add moves it).
#         Merge "to describe each variable..." text (bookmark removal already handled); now insert new
#         question paragraph after the "...etc.)" paragraph, before the existing blank numId=0 paragraph.
$rng16 = $d.Content
$rng16.Find.Execute("X2: own SAT /100, etc.)")
$para16 = $rng16.Paragraphs(1)
$prange16 = $para16.Range
$prange16.Collapse(0)
$prange16.InsertParagraphAfter()

$rng16b = $d.Content
$rng16b.Find.Execute("X2: own SAT /100, etc.)")
$newPara16 = $rng16b.Paragraphs(1).Next()
$newRange16 = $newPara16.Range
$newRange16.InsertAfter("Write down the regression equations for columns (3) and (6) in Table 2.5 of MM. Explain what is the purpose of the regressions presented in this table. ")

# 17. "on OVB and All things Regression by Friday Next Week]" -> "on CEF, OVB and All things Regression by Friday Next Week]"
ReplaceText "on OVB and All things Regression by Friday Next Week]" "on CEF, OVB and All things Regression by Friday Next Week]"

Write-Host "All edits applied"
